$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1013.2222
$ws.Range("I2").Value = 1199.6666
$ws.Range("J2").Value = 920
$ws.Range("K2").Value = 1199.6666
$ws.Range("L2").Value = 920
$ws.Range("M2").Value = -1086.6666
$ws.Range("N2").Value = -1146
$ws.Range("H28").Value = 284.86667
$ws.Range("I28").Value = 168
$ws.Range("K28").Value = 168
$ws.Range("M28").Value = 317
$ws.Range("H51").Value = 3549.8333
$ws.Range("I51").Value = 2759.8
$ws.Range("J51").Value = 7500
$ws.Range("K51").Value = 2759.8
$ws.Range("L51").Value = 7500
$ws.Range("M51").Value = -2275.8
$ws.Range("N51").Value = -8468
$ws.Range("H58").Value = 1504.1428
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300
$ws.Range("H106").Value = 17666.428
$ws.Range("I106").Value = 17666.428
$ws.Range("K106").Value = 17666.428
$ws.Range("M106").Value = -17035.428
$ws.Range("H107").Value = 456.86206
$ws.Range("I107").Value = 437.89285
$ws.Range("J107").Value = 988
$ws.Range("K107").Value = 437.89285
$ws.Range("L107").Value = 988
$ws.Range("M107").Value = 1482.10715
$ws.Range("N107").Value = -4828
$ws.Range("H113").Value = 6749.25
$ws.Range("I113").Value = 5666.3335
$ws.Range("J113").Value = 9998
$ws.Range("K113").Value = 5666.3335
$ws.Range("L113").Value = 9998
$ws.Range("M113").Value = -2412.3335
$ws.Range("N113").Value = -16506
$ws.Range("H116").Value = 9000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 9000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 9000
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -15884
$ws.Range("H138").Value = 12234.488
$ws.Range("J138").Value = 12594.182
$ws.Range("L138").Value = 37782.546
$ws.Range("N138").Value = -48062.546

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28771.107
$ws.Range("I32").Value = 26738.154
$ws.Range("J32").Value = 30533
$ws.Range("K32").Value = 26738.154
$ws.Range("L32").Value = 30533
$ws.Range("M32").Value = -26451.154
$ws.Range("N32").Value = -31107
$ws.Range("H36").Value = 13999.857
$ws.Range("I36").Value = 17600
$ws.Range("J36").Value = 4999.5
$ws.Range("K36").Value = 17600
$ws.Range("L36").Value = 4999.5
$ws.Range("M36").Value = -17254
$ws.Range("N36").Value = -5691.5
$ws.Range("H45").Value = 2552.75
$ws.Range("I45").Value = 2552.75
$ws.Range("K45").Value = 2552.75
$ws.Range("M45").Value = -2175.75
$ws.Range("H54").Value = 40000
$ws.Range("J54").Value = 40000
$ws.Range("L54").Value = 40000
$ws.Range("N54").Value = -41538
$ws.Range("H63").Value = 6800
$ws.Range("I63").Value = 5880
$ws.Range("K63").Value = 5880
$ws.Range("M63").Value = -5194
$ws.Range("H66").Value = 6800
$ws.Range("I66").Value = 5880
$ws.Range("K66").Value = 29400
$ws.Range("M66").Value = -25968
$ws.Range("H122").Value = 4033.6667
$ws.Range("I122").Value = 3081.1428
$ws.Range("J122").Value = 4867.125
$ws.Range("K122").Value = 9243.428400000001
$ws.Range("L122").Value = 14601.375
$ws.Range("M122").Value = -6793.428400000001
$ws.Range("N122").Value = -19501.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 669.73334
$ws.Range("I22").Value = 668
$ws.Range("J22").Value = 674.5
$ws.Range("K22").Value = 668
$ws.Range("L22").Value = 674.5
$ws.Range("M22").Value = -495
$ws.Range("N22").Value = -1020.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6250058
$ws.Range("I6").Value = 8333243.5
$ws.Range("K6").Value = 8333243.5
$ws.Range("M6").Value = -8333130.5
$ws.Range("H7").Value = 284.22223
$ws.Range("I7").Value = 315.6
$ws.Range("J7").Value = 245
$ws.Range("K7").Value = 315.6
$ws.Range("L7").Value = 245
$ws.Range("M7").Value = -202.6
$ws.Range("N7").Value = -471
$ws.Range("H22").Value = 332.46155
$ws.Range("I22").Value = 314.8889
$ws.Range("J22").Value = 372
$ws.Range("K22").Value = 314.8889
$ws.Range("L22").Value = 372
$ws.Range("M22").Value = 35.11110000000002
$ws.Range("N22").Value = -1072
$ws.Range("H31").Value = 6968.2856
$ws.Range("I31").Value = 6966.5
$ws.Range("K31").Value = 6966.5
$ws.Range("M31").Value = -6671.5
$ws.Range("H34").Value = 6968.2856
$ws.Range("I34").Value = 6966.5
$ws.Range("K34").Value = 6966.5
$ws.Range("M34").Value = -6764.5
$ws.Range("H58").Value = 3956
$ws.Range("J58").Value = 7589.3335
$ws.Range("L58").Value = 7589.3335
$ws.Range("N58").Value = -7995.3335
$ws.Range("H86").Value = 5944.636
$ws.Range("I86").Value = 3456.1738
$ws.Range("K86").Value = 3456.1738
$ws.Range("M86").Value = -2333.1738
$ws.Range("H89").Value = 5944.636
$ws.Range("I89").Value = 3456.1738
$ws.Range("K89").Value = 17280.869
$ws.Range("M89").Value = -11664.869
$ws.Range("H122").Value = 1782.7142
$ws.Range("I122").Value = 2293
$ws.Range("K122").Value = 6879
$ws.Range("M122").Value = -4429
$ws.Range("H136").Value = 3956
$ws.Range("J136").Value = 7589.3335
$ws.Range("L136").Value = 22768.0005
$ws.Range("N136").Value = -27868.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 748.5
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H56").Value = 9962.647000000001
$ws.Range("I56").Value = 9962.647000000001
$ws.Range("K56").Value = 9962.647000000001
$ws.Range("M56").Value = -9432.647000000001
$ws.Range("H98").Value = 312.77777
$ws.Range("I98").Value = 390.5
$ws.Range("J98").Value = 250.6
$ws.Range("K98").Value = 1171.5
$ws.Range("L98").Value = 751.8
$ws.Range("M98").Value = 326.5
$ws.Range("N98").Value = -3747.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1482.3549
$ws.Range("I102").Value = 749.3200000000001
$ws.Range("J102").Value = 4536.6665
$ws.Range("K102").Value = 749.3200000000001
$ws.Range("L102").Value = 4536.6665
$ws.Range("M102").Value = 872.6799999999999
$ws.Range("N102").Value = -7780.6665
$ws.Range("H136").Value = 50531.89
$ws.Range("J136").Value = 50531.89
$ws.Range("L136").Value = 151595.67
$ws.Range("N136").Value = -156695.67

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2805.3
$ws.Range("I46").Value = 2167.6
$ws.Range("J46").Value = 3443
$ws.Range("K46").Value = 2167.6
$ws.Range("L46").Value = 3443
$ws.Range("M46").Value = -1979.6
$ws.Range("N46").Value = -3819
$ws.Range("H82").Value = 1912.1765
$ws.Range("I82").Value = 2060.6667
$ws.Range("J82").Value = 1555.8
$ws.Range("K82").Value = 2060.6667
$ws.Range("L82").Value = 1555.8
$ws.Range("M82").Value = -1699.6667
$ws.Range("N82").Value = -2277.8
$ws.Range("H85").Value = 1912.1765
$ws.Range("I85").Value = 2060.6667
$ws.Range("J85").Value = 1555.8
$ws.Range("K85").Value = 2060.6667
$ws.Range("L85").Value = 1555.8
$ws.Range("M85").Value = -812.6667000000002
$ws.Range("N85").Value = -4051.8
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6520.684
$ws.Range("I62").Value = 5998.5
$ws.Range("K62").Value = 5998.5
$ws.Range("M62").Value = -5374.5
$ws.Range("H65").Value = 6520.684
$ws.Range("I65").Value = 5998.5
$ws.Range("K65").Value = 29992.5
$ws.Range("M65").Value = -26872.5
$ws.Range("H81").Value = 4444.222
$ws.Range("I81").Value = 3999.875
$ws.Range("K81").Value = 7999.75
$ws.Range("M81").Value = -6938.75
$ws.Range("H84").Value = 4444.222
$ws.Range("I84").Value = 3999.875
$ws.Range("K84").Value = 39998.75
$ws.Range("M84").Value = -34694.75
